# Apply updated taxon-observation data to rows 5-7 of the active sheet.
# The edit rotates the species data (columns E,F,G,H,Q,R) and the Id/
# Taxonsorteringsordning columns (A,B) between rows 5, 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - becomes the former "Västlig hakmossa" record
$ws.Range("A5").Value = 112395251
$ws.Range("B5").Value = 93230
$ws.Range("E5").Value = 2810
$ws.Range("F5").Value = "Västlig hakmossa"
$ws.Range("G5").Value = "Rhytidiadelphus loreus"
$ws.Range("H5").Value = "(Hedw.) Warnst."
$ws.Range("Q5").Value = 333055
$ws.Range("R5").Value = 6626785

# Row 6 - becomes the former "Vedticka" record
$ws.Range("A6").Value = 112395254
$ws.Range("B6").Value = 89517
$ws.Range("E6").Value = 5447
$ws.Range("F6").Value = "Vedticka"
$ws.Range("G6").Value = "Fuscoporia viticola"
$ws.Range("H6").Value = "(Schwein.) Murrill"
$ws.Range("Q6").Value = 333021
$ws.Range("R6").Value = 6626691

# Row 7 - becomes the former "Korallblylav" record
$ws.Range("A7").Value = 112395253
$ws.Range("B7").Value = 78671
$ws.Range("E7").Value = 229497
$ws.Range("F7").Value = "Korallblylav"
$ws.Range("G7").Value = "Parmeliella triptophylla"
$ws.Range("H7").Value = "(Ach.) Müll.Arg."
$ws.Range("Q7").Value = 333022
$ws.Range("R7").Value = 6626740
